# Fruta / hortaliza, semanal
# Weekly update: a new price-record row is inserted at row 19 (pushing the
# existing rows 19-37 down to 20-38) with the latest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 19; existing rows 19:37 shift down to 20:38.
$ws.Rows("19:19").Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A19").Value = 7
$ws.Range("B19").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C19").Value = "Ñuble"
$ws.Range("D19").Value = 44540
$ws.Range("E19").Value = 16
$ws.Range("F19").Value = 100112026
$ws.Range("G19").Value = "Haba"
$ws.Range("H19").Value = "Sin especificar"
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 100
$ws.Range("K19").Value = 6500
$ws.Range("L19").Value = 7000
$ws.Range("M19").Value = 6750
$ws.Range("N19").Value = "$/saco 25 kilos"
$ws.Range("O19").Value = "Provincia de Diguillín"
$ws.Range("P19").Value = 270
$ws.Range("Q19").Value = 25
$ws.Range("R19").Value = "Hortaliza"
